$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value changes (row permutation + content updates) ---
$ws.Range("A2").Value = 111487420
$ws.Range("Q2").Value = 626221.6215819545
$ws.Range("R2").Value = 6893253.628042356
$ws.Range("A3").Value = 111487418
$ws.Range("B3").Value = 77515
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = 'Garnlav'
$ws.Range("G3").Value = 'Alectoria sarmentosa'
$ws.Range("H3").Value = '(Ach.) Ach.'
$ws.Range("Q3").Value = 626244.4413132126
$ws.Range("R3").Value = 6893219.854707362
$ws.Range("AC3").Value = ""
$ws.Range("A4").Value = 111487429
$ws.Range("B4").Value = 96348
$ws.Range("D4").Value = 'VU'
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = 'Knärot'
$ws.Range("G4").Value = 'Goodyera repens'
$ws.Range("H4").Value = '(L.) R. Br.'
$ws.Range("Q4").Value = 626289.1696174983
$ws.Range("R4").Value = 6893239.309586792
$ws.Range("AC4").Value = 'ca 50 bladrosetter inom 2m2'
$ws.Range("A5").Value = 111487423
$ws.Range("B5").Value = 73696
$ws.Range("E5").Value = 6440
$ws.Range("F5").Value = 'Vitgrynig nållav'
$ws.Range("G5").Value = 'Chaenotheca subroscida'
$ws.Range("H5").Value = '(Eitner) Zahlbr.'
$ws.Range("Q5").Value = 626222.1012433186
$ws.Range("R5").Value = 6893215.74357231
$ws.Range("A6").Value = 111487415
$ws.Range("B6").Value = 73634
$ws.Range("E6").Value = 6426
$ws.Range("F6").Value = 'Kattfotslav'
$ws.Range("G6").Value = 'Felipes leucopellaeus'
$ws.Range("H6").Value = '(Ach.) Frisch & G.Thor'
$ws.Range("Q6").Value = 626321.4062460049
$ws.Range("R6").Value = 6893191.850843907
$ws.Range("A7").Value = 111487428
$ws.Range("B7").Value = 96348
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = 'Knärot'
$ws.Range("G7").Value = 'Goodyera repens'
$ws.Range("H7").Value = '(L.) R. Br.'
$ws.Range("Q7").Value = 626201.261590388
$ws.Range("R7").Value = 6893121.376245681
$ws.Range("AC7").Value = 'ca 50 bladrosetter spridda över 2 m2'
$ws.Range("A8").Value = 111487419
$ws.Range("B8").Value = 5135
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 105930
$ws.Range("F8").Value = 'Vågbandad barkbock'
$ws.Range("G8").Value = 'Semanotus undatus'
$ws.Range("H8").Value = '(Linnaeus, 1758)'
$ws.Range("Q8").Value = 626208.9904600172
$ws.Range("R8").Value = 6893240.521565447
$ws.Range("A9").Value = 111487427
$ws.Range("B9").Value = 78578
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 6458
$ws.Range("F9").Value = 'Lunglav'
$ws.Range("G9").Value = 'Lobaria pulmonaria'
$ws.Range("H9").Value = '(L.) Hoffm.'
$ws.Range("Q9").Value = 626206.7588566126
$ws.Range("R9").Value = 6893112.222867905
$ws.Range("A10").Value = 111487422
$ws.Range("B10").Value = 89369
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 5447
$ws.Range("F10").Value = 'Vedticka'
$ws.Range("G10").Value = 'Fuscoporia viticola'
$ws.Range("H10").Value = '(Schwein.) Murrill'
$ws.Range("Q10").Value = 626214.6507017991
$ws.Range("R10").Value = 6893264.597594698
$ws.Range("A11").Value = 111487416
$ws.Range("B11").Value = 96348
$ws.Range("D11").Value = 'VU'
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = 'Knärot'
$ws.Range("G11").Value = 'Goodyera repens'
$ws.Range("H11").Value = '(L.) R. Br.'
$ws.Range("Q11").Value = 626241.4078639001
$ws.Range("R11").Value = 6893275.892764967
$ws.Range("A12").Value = 111487421
$ws.Range("B12").Value = 89405
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 1202
$ws.Range("F12").Value = 'Ullticka'
$ws.Range("G12").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H12").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q12").Value = 626220.5943774013
$ws.Range("R12").Value = 6893268.563298941
$ws.Range("A13").Value = 111487424
$ws.Range("B13").Value = 96348
$ws.Range("D13").Value = 'VU'
$ws.Range("E13").Value = 220787
$ws.Range("F13").Value = 'Knärot'
$ws.Range("G13").Value = 'Goodyera repens'
$ws.Range("H13").Value = '(L.) R. Br.'
$ws.Range("Q13").Value = 626179.4687150358
$ws.Range("R13").Value = 6893140.215102527
$ws.Range("A14").Value = 111487417
$ws.Range("Q14").Value = 626274.1114414346
$ws.Range("R14").Value = 6893228.451636742
$ws.Range("AC14").Value = ""

# --- AF column placeholder cells follow the same row permutation. ---
# AF14 keeps its (empty) placeholder cell before and after the edit, so use
# it as a template to stamp empty placeholder cells onto AF4, AF7 and AF11,
# then blank out the placeholders that should no longer be present.
$ws.Range("AF14").Copy($ws.Range("AF4"))
$ws.Range("AF14").Copy($ws.Range("AF7"))
$ws.Range("AF14").Copy($ws.Range("AF11"))
$ws.Range("AF3").Value = ""
$ws.Range("AF8").Value = ""
$ws.Range("AF9").Value = ""
